{"js": "// Find the paragraph that currently starts with \"August 6, 2024 Dear\" and:\n//  1) insert a new paragraph \"June 28, 2024\" right before it (placeholder systemDate,\n//     formatted from YYYY-MM-DD to \"Month DD, YYYY\"); and\n//  2) strip the old \"August 6, 2024 \" date prefix from that paragraph so it just\n//     starts with \"Dear\".\nconst results = context.document.body.search(\"August 6, 2024 Dear\", { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const hit = results.items[0];\n  const targetPara = hit.paragraphs.getFirst();\n\n  // Mimics: systemDate \"2024-06-28\" run through a YYYY-MM-DD -> \"Month DD, YYYY\"\n  // formatter to produce formattedSystemDate.\n  const systemDate = \"2024-06-28\";\n  const [year, month, day] = systemDate.split(\"-\").map((p) => parseInt(p, 10));\n  const monthNames = [\n    \"January\", \"February\", \"March\", \"April\", \"May\", \"June\",\n    \"July\", \"August\", \"September\", \"October\", \"November\", \"December\"\n  ];\n  const formattedSystemDate = `${monthNames[month - 1]} ${day}, ${year}`;\n\n  const newPara = targetPara.insertParagraph(formattedSystemDate, Word.InsertLocation.before);\n  newPara.font.name = \"Century Gothic\";\n  newPara.paragraphFormat.rightIndent = 14.35; // 287 twips\n\n  // Replace \"August 6, 2024 Dear\" with just \"Dear\" (keeps the rest of the\n  // paragraph - \" JAHIR,\" - untouched).\n  hit.insertText(\"Dear\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Placeholder systemDate (YYYY-MM-DD) and a small helper that formats it as\n# \"Month DD, YYYY\" -> formattedSystemDate.\n$systemDate = \"2024-06-28\"\nfunction Format-SystemDate($isoDate) {\n    $parts = $isoDate.Split(\"-\")\n    $year = [int]$parts[0]\n    $month = [int]$parts[1]\n    $day = [int]$parts[2]\n    $monthNames = @(\"January\", \"February\", \"March\", \"April\", \"May\", \"June\", `\n                     \"July\", \"August\", \"September\", \"October\", \"November\", \"December\")\n    $monthName = $monthNames[$month - 1]\n    return \"$monthName $day, $year\"\n}\n$formattedSystemDate = Format-SystemDate $systemDate\n\n# Locate the paragraph that currently reads \"August 6, 2024 Dear JAHIR,\".\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"August 6, 2024 Dear\")\n\nif ($found) {\n    $targetParagraph = $findRange.Paragraphs(1)\n    $insertPoint = $targetParagraph.Range.Duplicate\n    $insertPoint.Collapse(1)  # wdCollapseStart\n\n    # Insert a new paragraph with the formatted system date right before the\n    # existing date/greeting paragraph.\n    $insertPoint.InsertBefore($formattedSystemDate + \"`r\")\n    $insertPoint.ParagraphFormat.RightIndent = 14.35  # 287 twips\n    $insertPoint.Font.Name = \"Century Gothic\"\n\n    # Strip the old \"August 6, 2024 \" date prefix, leaving just \"Dear\" (the\n    # rest of the paragraph - \" JAHIR,\" - is untouched).\n    $replaceRange = $d.Content\n    $replaceRange.Find.Execute(\"August 6, 2024 Dear\") | Out-Null\n    $replaceRange.Text = \"Dear\"\n}\n"}
